$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.491.99"
$ws.Range("E2").Value = "  -1.19%  "

$ws.Range("D3").Value = "1.911.75"
$ws.Range("E3").Value = "  -1.49%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.62"
$ws.Range("E5").Value = "  -1.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4778"
$ws.Range("E7").Value = "  -2.71%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2846"
$ws.Range("E8").Value = "  -3.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06710"
$ws.Range("E9").Value = "  -2.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.55"
$ws.Range("E10").Value = "  +1.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "103.39"
$ws.Range("E11").Value = "  -2.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07757"
$ws.Range("E12").Value = "  -0.12%  "

$ws.Range("D13").Value = "1.915.35"
$ws.Range("E13").Value = "  -1.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.193"
$ws.Range("E14").Value = "  -3.28%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6697"
$ws.Range("E15").Value = "  -4.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "276.36"
$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("D17").Value = "30.487.63"
$ws.Range("E17").Value = "  -1.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9993"
$ws.Range("E18").Value = "  -0.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007492"
$ws.Range("E19").Value = "  -3.06%  "

$ws.Range("E20").Value = "  -3.78%  "

$ws.Range("E21").Value = "  -4.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9988"
$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.294"
$ws.Range("E23").Value = "  -3.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.360"
$ws.Range("E24").Value = "  -4.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.70"
$ws.Range("E25").Value = "  +0.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.22"
$ws.Range("E26").Value = "  -1.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.082"
$ws.Range("E27").Value = "  -3.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.385"
$ws.Range("E28").Value = "  -0.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09983"
$ws.Range("E29").Value = "  -4.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.595"
$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.512"
$ws.Range("E31").Value = "  -3.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.256"
$ws.Range("E32").Value = "  -2.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04711"
$ws.Range("E33").Value = "  -3.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7272"
$ws.Range("E34").Value = "  -3.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.117"
$ws.Range("E35").Value = "  -3.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.718"
$ws.Range("E36").Value = "  -0.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01907"
$ws.Range("E37").Value = "  -4.71%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.609"
$ws.Range("E38").Value = "  -1.79%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.371"
$ws.Range("E39").Value = "  -1.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "74.40"
$ws.Range("E40").Value = "  -5.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.957"
$ws.Range("E41").Value = "  -6.50%  "

$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.60"
$ws.Range("E42").Value = "  -1.03%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8613"
$ws.Range("E43").Value = "  -5.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4262"
$ws.Range("E44").Value = "  -3.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9991"
$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.411"
$ws.Range("E46").Value = "  -3.50%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "950.45"
$ws.Range("E47").Value = "  -4.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1206"
$ws.Range("E48").Value = "  -3.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.63"
$ws.Range("E49").Value = "  -3.99%  "

$ws.Range("E50").Value = "  +0.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.727"
$ws.Range("E51").Value = "  -5.07%  "
